# "add solo titles dungeon bomb balloon sort"
# Append 3 new player/title rows to the bottom of the existing data table
# on Sheet1 (rows 83-85), then move the selection to reflect where the
# author ended up working (G15), matching the scrolled/selected state
# captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows: player, ratelaw, ratehigh, Accuracyrate, time, team1..team5
$newRows = @(
    @("無法者集団", 2400, 9999, 1, 10, 10001, 944, 270, 996, 948),
    @("最高戦力",   2400, 9999, 1, 10,   830, 357, 595, 538, 1678),
    @("知識王",     2400, 9999, 1, 10,  1076,   2,   1,   3,   17)
)

$startRow = 83
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value() = $rowData[$c - 1]
    }
}

# Leave the selection on the cell the author last touched.
$ws.Range("G15").Select()
